$wb = $excel.ActiveWorkbook

# zh-cn sheet: update the handoff/handback datetime pair for the 60e486bb... row
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D4").Value = "2016-02-22 04:28:35"
$wsZh.Range("G4").Value = "2016-02-22 04:29:40"

# de-de sheet: update the handoff/handback datetime pair for the 60e486bb... row
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D4").Value = "2016-02-22 04:28:50"
$wsDe.Range("G4").Value = "2016-02-22 04:30:05"
